# Insert a new weekly record row above current row 198 (Primera / Región del
# Maule, fecha 44617) and push the existing rows 198-218 down to 199-219,
# keeping all of their original data intact (this is how the source system
# adds a freshly scraped day to the middle of the historical series).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("198:198").Insert()

$ws.Range("A198").Value = 11
$ws.Range("B198").Value = "Vega Monumental Concepción"
$ws.Range("C198").Value = "Bíobío"
$ws.Range("D198").Value = 44617
$ws.Range("E198").Value = 8
$ws.Range("F198").Value = 100112008
$ws.Range("G198").Value = "Coliflor"
$ws.Range("H198").Value = "Sin especificar"
$ws.Range("I198").Value = "Primera"
$ws.Range("J198").Value = 3000
$ws.Range("K198").Value = 700
$ws.Range("L198").Value = 800
$ws.Range("M198").Value = 767
$ws.Range("N198").Value = "`$/unidad"
$ws.Range("O198").Value = "Región del Maule"
$ws.Range("P198").Value = 767
$ws.Range("Q198").Value = 1
$ws.Range("R198").Value = "Hortaliza"
